$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '88.327.49'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').Value = '3.023.72'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = "'208.27"
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').Value = "'606.12"
$ws.Range('E6').Value = '  -3.28%  '
$ws.Range('E7').Value = '  -6.93%  '
$ws.Range('D8').Value = "'0.872"
$ws.Range('E8').Value = '  +23.00%  '
$ws.Range('D9').Value = "'1.00"
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = '3.021.46'
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('E11').Value = '  +19.08%  '
$ws.Range('E12').Value = '  +3.63%  '
$ws.Range('E13').Value = '  -5.39%  '
$ws.Range('D14').Value = "'5.32"
$ws.Range('E14').Value = '  +2.37%  '
$ws.Range('D15').Value = '88.306.19'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '3.600.73'
$ws.Range('E16').Value = '  -1.37%  '
$ws.Range('D17').Value = "'31.41"
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '3.064.91'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').Value = "'3.34"
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').Value = "'0.0000202"
$ws.Range('E20').Value = '  -1.38%  '
$ws.Range('E21').Value = '  +2.00%  '
$ws.Range('D22').Value = "'418.92"
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('E23').Value = '  +2.50%  '
$ws.Range('D24').Value = "'7.97"
$ws.Range('E24').Value = '  -2.95%  '
$ws.Range('D25').Value = "'5.35"
$ws.Range('E25').Value = '  +3.79%  '
$ws.Range('E26').Value = '  +5.74%  '
$ws.Range('D27').Value = "'11.44"
$ws.Range('E27').Value = '  +1.94%  '
$ws.Range('D28').Value = '3.204.45'
$ws.Range('E28').Value = '  -1.66%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  +9.28%  '
$ws.Range('E31').Value = '  +2.24%  '
$ws.Range('D32').Value = "'8.08"
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').Value = "'496.18"
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('E34').Value = '  -9.52%  '
$ws.Range('E35').Value = '  -2.66%  '
$ws.Range('E36').Value = '  -1.51%  '
$ws.Range('D37').Value = "'22.29"
$ws.Range('E37').Value = '  +3.66%  '
$ws.Range('E38').Value = '  -1.71%  '
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('E40').Value = '  +5.10%  '
$ws.Range('D41').Value = "'1.00"
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  +11.17%  '
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('D45').Value = "'145.67"
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('D47').Value = "'43.30"
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').Value = "'0.0678"
$ws.Range('E48').Value = '  +13.31%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = "'3.99"
$ws.Range('E49').Value = '  +3.12%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').Value = "'1.19"
$ws.Range('E50').Value = '  +2.76%  '
$ws.Range('D51').Value = "'154.30"
$ws.Range('E51').Value = '  -5.39%  '
